$d = $word.ActiveDocument

# The paragraph "Process of creating features..." currently ends with the
# hidden _GoBack bookmark glued right after its text. We need to split that
# single paragraph into three:
#   1) "Process of creating features that enhance the performance of ML models"
#   2) "- Modify preexisting features"   <- keeps the _GoBack bookmark
#   3) "- Design new features"

$bm = $d.Bookmarks.Item("_GoBack")

# Insert the first new line of text through the bookmark's own Range so the
# (collapsed) bookmark expands to wrap the freshly typed text, i.e. it ends
# up positioned right after "- Modify preexisting features".
$bm.Range.InsertAfter("- Modify preexisting features")

# Find the boundary between "...ML models" and "- Modify preexisting
# features" (they are currently glued together in the same run/paragraph)
# and insert a paragraph break there, splitting it into two paragraphs.
$marker = "ML models- Modify preexisting features"
$full = $d.Content.Text
$markerIdx = $full.IndexOf($marker)
$splitPos = $markerIdx + "ML models".Length
$d.Range($splitPos, $splitPos).InsertAfter("`r")

# Append the third paragraph after the (now relocated) bookmark, using a
# plain Range rather than the bookmark's own Range so the bookmark itself
# does not expand to swallow this new text.
$bm = $d.Bookmarks.Item("_GoBack")
$endPos = $bm.End
$d.Range($endPos, $endPos).InsertAfter("`r- Design new features")
